# "settel bahagian popup utk parliament"
# The workbook was re-saved from a newer Excel, which (re)computed "best fit"
# column widths for Sheet1 (based on the actual cell contents/fonts) and left
# the selection/active cell parked on O1 instead of the whole-sheet selection.
# Reproduce the two user-visible effects via the Excel object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Best-fit column widths for columns A:AG (1-33), expressed as the
# Range.ColumnWidth value that reproduces the stored worksheet width.
$colWidths = @{
    1  = 29.0221354166667
    2  = 27.0221354166667
    3  = 5.16666666666667
    4  = 5.16666666666667
    5  = 5.16666666666667
    6  = 5.16666666666667
    7  = 5.16666666666667
    8  = 5.16666666666667
    9  = 5.16666666666667
    10 = 5.16666666666667
    11 = 4.73697916666667
    12 = 5.16666666666667
    13 = 4.73697916666667
    14 = 5.16666666666667
    15 = 4.73697916666667
    16 = 4.16666666666667
    17 = 3.30729166666667
    18 = 4.16666666666667
    19 = 8.73697916666667
    20 = 5.16666666666667
    21 = 10.4518229166667
    22 = 7.16666666666667
    23 = 16.4518229166667
    24 = 8.87760416666667
    25 = 10.8776041666667
    26 = 17.7369791666667
    27 = 13.3072916666667
    28 = 8.59244791666667
    29 = 12.8776041666667
    30 = 8.59244791666667
    31 = 12.8776041666667
    32 = 15.0221354166667
    33 = 7.45182291666667
}

foreach ($col in $colWidths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $colWidths[$col]
}

# Move the active cell / selection from the whole data range to O1.
$ws.Range("O1").Select() | Out-Null
